$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.154.70'
$ws.Range('E2').Value = '  +2.58%  '
$ws.Range('D3').Value = '1.912.41'
$ws.Range('E3').Value = '  +2.22%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9985'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '333.73'
$ws.Range('E5').Value = '  -1.58%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9990'
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4645'
$ws.Range('E7').Value = '  -1.24%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.4107'
$ws.Range('E8').Value = '  +3.39%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '47.73'
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.08025'
$ws.Range('E10').Value = '  -0.11%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.011'
$ws.Range('E11').Value = '  +1.00%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '21.91'
$ws.Range('E12').Value = '  -0.43%  '
$ws.Range('D13').Value = '1.896.35'
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.960'
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.112'
$ws.Range('E15').Value = '  -2.13%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '89.31'
$ws.Range('E16').Value = '  -2.08%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.9984'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.00001035'
$ws.Range('E18').Value = '  -0.49%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06587'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.58'
$ws.Range('E20').Value = '  -0.16%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.9998'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').Value = '29.174.35'
$ws.Range('E22').Value = '  +2.65%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.452'
$ws.Range('E23').Value = '  -0.47%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.34'
$ws.Range('E24').Value = '  +2.40%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.236'
$ws.Range('E25').Value = '  -0.75%  '
$ws.Range('D26').Value = '2.124.80'
$ws.Range('E26').Value = '  +1.56%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '157.30'
$ws.Range('E27').Value = '  -2.09%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.78'
$ws.Range('E28').Value = '  +0.07%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.121'
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '5.440'
$ws.Range('E30').Value = '  -1.26%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '118.47'
$ws.Range('E31').Value = '  -1.60%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.9906'
$ws.Range('E32').Value = '  +1.60%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.09426'
$ws.Range('E33').Value = '  -0.85%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.438'
$ws.Range('E34').Value = '  +4.49%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.590'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.330'
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.06110'
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02252'
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '8.400'
$ws.Range('E39').Value = '  +0.35%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.180'
$ws.Range('E40').Value = '  +0.13%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.5831'
$ws.Range('E41').Value = '  -2.21%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.9985'
$ws.Range('E42').Value = '  -0.16%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '10.23'
$ws.Range('E43').Value = '  -1.16%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.1831'
$ws.Range('E44').Value = '  -2.67%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.275'
$ws.Range('E45').Value = '  -0.67%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.357'
$ws.Range('E46').Value = '  +14.01%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5527'
$ws.Range('E47').Value = '  -1.26%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '12.12'
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.925'
$ws.Range('E49').Value = '  -1.67%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.07059'
$ws.Range('E50').Value = '  +1.79%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '47.95'
$ws.Range('E51').Value = '  +22.71%  '
